$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.911.41'
$ws.Range("E2").Value = '  +1.96%  '
$ws.Range("D3").Value = '3.632.68'
$ws.Range("E3").Value = '  +3.74%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''606.95'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.74%  '
$ws.Range("D6").Value = '''199.83'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("D7").Value = '''0.628'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +1.21%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  +9.84%  '
$ws.Range("D10").Value = '''0.648'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.57%  '
$ws.Range("D11").Value = '''53.75'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.41%  '
$ws.Range("D12").Value = '''0.0000306'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.11%  '
$ws.Range("D13").Value = '''9.57'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").Value = '4.199.44'
$ws.Range("E14").Value = '  +3.45%  '
$ws.Range("D15").Value = '''680.87'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +14.73%  '
$ws.Range("D16").Value = '''12.98'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.67%  '
$ws.Range("D17").Value = '70.913.08'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").Value = '3.634.63'
$ws.Range("E18").Value = '  +3.90%  '
$ws.Range("D19").Value = '''19.06'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("E20").Value = '  +0.13%  '
$ws.Range("D21").Value = '''1.00'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.03%  '
$ws.Range("D22").Value = '''18.67'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.68%  '
$ws.Range("D23").Value = '''5.40'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +3.31%  '
$ws.Range("D24").Value = '''105.86'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +4.09%  '
$ws.Range("D25").Value = '''4.63'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").Value = '''3.01'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.97%  '
$ws.Range("D27").Value = '''10.52'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("D28").Value = '''9.90'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.77%  '
$ws.Range("D29").Value = '''34.39'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +4.49%  '
$ws.Range("D30").Value = '''4.63'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +8.63%  '
$ws.Range("D31").Value = '''7.20'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.25%  '
$ws.Range("D32").Value = '''12.21'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.58%  '
$ws.Range("D33").Value = '''0.116'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.29%  '
$ws.Range("D34").Value = '''63.31'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("B35").Value = 'PEPE'
$ws.Range("C35").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D35").Value = '0.0₃0873'
$ws.Range("E35").Value = '  +7.58%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '3.953.70'
$ws.Range("E36").Value = '  +5.96%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").Value = '''3.02'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.69%  '
$ws.Range("D39").Value = '''36.71'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '''503.73'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.07%  '
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").Value = '''0.388'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '''3.55'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.06%  '
$ws.Range("E43").Value = '  +2.93%  '
$ws.Range("E44").Value = '  +9.67%  '
$ws.Range("D45").Value = '''0.0459'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +2.37%  '
$ws.Range("D46").Value = '''3.50'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +6.62%  '
$ws.Range("D47").Value = '''0.141'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.60%  '
$ws.Range("D48").Value = '''8.66'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.32%  '
$ws.Range("E49").Value = '  -0.43%  '
$ws.Range("D50").Value = '''0.000248'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.85%  '
$ws.Range("D51").Value = '''1.30'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.04%  '
